# Remove the initial "New" product-status sub-bullet and its explanatory
# sub-sub-bullet ("The product was created within the last 7 days") that
# immediately follow the "Product Status" bullet, leaving "In Stock" (and
# its explanation) as the first status entry.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($text -eq "New") {
        $startPara = $p
    } elseif ($text -eq "The product was created within the last 7 days") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
